$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New string values are introduced to the shared-string table in the order
# the cells were authored, so write the cells in that exact order.
$ws.Range("A12").Value = "Duckweed"
$ws.Range("H12").Value = "Small ponds"
$ws.Range("I12").Value = "water bodies"
$ws.Range("A13").Value = "Vetiver Grass"
$ws.Range("B13").Value = "Chrysopogon zizanioides"
$ws.Range("B12").Value = "Lemna minor"
$ws.Range("H13").Value = "Sandy/Loamy,Riverbanks"
$ws.Range("I13").Value = "erosion-prone areas"
$ws.Range("E12").Value = "Agricultural canals"
$ws.Range("E13").Value = "Roadsides and erosion-prone areas"

# Numeric cells
$ws.Range("C12").Value = 10
$ws.Range("D12").Value = 0.1
$ws.Range("F12").Value = 0.85
$ws.Range("G12").Value = 1.3

$ws.Range("C13").Value = 15
$ws.Range("D13").Value = 20
$ws.Range("F13").Value = 0.9
$ws.Range("G13").Value = 1.2

# Update the active cell selection to match the post-edit state
$ws.Range("E13").Select() | Out-Null
